$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row -> source (old) row, derived from the permutation described by the diff
$map = @{}
$map[2] = 16
$map[3] = 17
$map[4] = 132
$map[5] = 41
$map[6] = 94
$map[7] = 33
$map[8] = 88
$map[9] = 25
$map[10] = 72
$map[11] = 48
$map[12] = 113
$map[13] = 20
$map[14] = 103
$map[15] = 14
$map[16] = 58
$map[17] = 78
$map[18] = 70
$map[19] = 96
$map[20] = 2
$map[21] = 97
$map[22] = 87
$map[23] = 24
$map[24] = 39
$map[25] = 64
$map[26] = 107
$map[27] = 115
$map[28] = 102
$map[29] = 30
$map[30] = 89
$map[31] = 76
$map[32] = 42
$map[33] = 90
$map[34] = 82
$map[35] = 116
$map[36] = 68
$map[37] = 6
$map[38] = 27
$map[39] = 134
$map[40] = 110
$map[41] = 51
$map[42] = 124
$map[43] = 104
$map[44] = 80
$map[45] = 31
$map[46] = 133
$map[47] = 45
$map[48] = 63
$map[49] = 19
$map[50] = 100
$map[51] = 71
$map[52] = 12
$map[53] = 123
$map[54] = 98
$map[55] = 125
$map[56] = 83
$map[57] = 53
$map[58] = 121
$map[59] = 109
$map[60] = 29
$map[61] = 11
$map[62] = 44
$map[63] = 55
$map[64] = 46
$map[65] = 73
$map[66] = 106
$map[67] = 101
$map[68] = 99
$map[69] = 37
$map[70] = 60
$map[71] = 81
$map[72] = 129
$map[73] = 105
$map[74] = 28
$map[75] = 93
$map[76] = 59
$map[77] = 7
$map[78] = 128
$map[79] = 108
$map[80] = 122
$map[81] = 85
$map[82] = 50
$map[83] = 74
$map[84] = 118
$map[85] = 8
$map[86] = 52
$map[87] = 62
$map[88] = 38
$map[89] = 130
$map[90] = 111
$map[91] = 69
$map[92] = 35
$map[93] = 43
$map[94] = 54
$map[95] = 117
$map[96] = 114
$map[97] = 126
$map[98] = 65
$map[99] = 22
$map[100] = 131
$map[101] = 66
$map[102] = 36
$map[103] = 4
$map[104] = 15
$map[105] = 112
$map[106] = 13
$map[107] = 61
$map[108] = 127
$map[109] = 84
$map[110] = 34
$map[111] = 75
$map[112] = 57
$map[113] = 120
$map[114] = 49
$map[115] = 77
$map[116] = 91
$map[117] = 3
$map[118] = 95
$map[119] = 9
$map[120] = 26
$map[121] = 18
$map[122] = 79
$map[123] = 92
$map[124] = 47
$map[125] = 32
$map[126] = 67
$map[127] = 10
$map[128] = 5
$map[129] = 86
$map[130] = 21
$map[131] = 119
$map[132] = 40
$map[133] = 56
$map[134] = 23

# Columns that carry per-row data values (others - A,B,C,E,F,G,H,I,R - are constant across rows)
$cols = @(4, 10, 11, 12, 13, 14, 15, 16, 17)  # D, J, K, L, M, N, O, P, Q

# Snapshot every relevant cell BEFORE any writes, so source data is not clobbered mid-pass
$snapshot = @{}
for ($r = 2; $r -le 134; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowData
}

# Apply the permutation: row $r gets the snapshot of row $map[$r]
for ($r = 2; $r -le 134; $r++) {
    $srcRow = $map[$r]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value2 = $srcData[$c]
    }
}

"Done applying permutation to rows 2-134"